$d = $word.ActiveDocument

$pairs = @(
    @("教育是人类生活中最关键的方面之一。它是一个工具，可以帮助人们获得知识、技能、价值观和态度，这些都是他们个人和职业发展所必需的。教育不仅限于学校和大学，它可以发生在任何地方、任何时候。它是一个不断进行的过程，可以帮助个人发展潜力、改善生活质量并为社会做出贡献。", "教育是人类生命中最关键的方面之一。它是一种工具，帮助人们获取知识、技能、价值观及态度，这些都是他们个人和职业发展所必需的。教育并不局限于学校和大学，它可以在任何地方和任何时候发生。它是一个不断的过程，帮助个人发展他们的潜力，改善他们的生活质量，并为社会做出贡献。"),
    @("教育对于个人的成长和发展是至关重要的。它帮助个人更好地了解自己，了解自己的优势、劣势、兴趣和价值观。它使个人能够设定目标，做出明智的决策，并对自己的行为负责。教育还可以帮助个人发展批判性思维能力、创造力和解决问题的能力。这些能力在今天的世界中至关重要，因为适应变化和解决复杂问题的能力是成功的关键。", "教育对个人成长和发展至关重要。它帮助人们更好地了解自己，了解自己的优势、劣势、兴趣和价值观。它使人们能够设定目标、做出明智的决策并对自己的行为负责。教育还有助于个人发展批判性思维技能、创造力和解决问题的能力。这些技能对于在今天的世界中成功至关重要，因为适应变化和解决复杂问题的能力是至关重要的。"),
    @("教育在社会和经济发展中也起着重要作用。它帮助个人获得就业和创业所必需的技能和知识。教育对于创造一个有技能和胜任力的劳动力来推动经济增长和发展而言至关重要。它也对于社会发展至关重要，因为它帮助个人理解和欣赏多样性，促进社会凝聚力，并培养社区意识。", "教育在社会和经济发展中也起着重要作用。它帮助个人获得就业和创业所必需的技能和知识。教育对于创造技能娴熟和胜任能力的劳动力至关重要，能为经济的增长和发展做出贡献。同时，教育对于社会发展同样至关重要，因为它有助于个人理解和欣赏多样性，促进社会凝聚力，并培养社区意识。"),
    @("此外，教育是一项基本权利，应该人人可及。它是实现公平公正社会的基石，让每个人都有平等机会成功。教育为个人提供改善生活所需的工具和资源，是降低贫困和不平等的有力手段。教育对于促进民主、人权和社会正义也至关重要。", "此外，教育是应该向每个人开放的基本权利。它是一个公平正义社会的基础，在这里每个人都有平等的成功机会。教育为个人提供了改善生活所需的工具和资源，同时它也是减少贫困和不平等的有力手段。教育对于促进民主、人权和社会公正也是至关重要的。"),
    @("总之，教育是人类生活中至关重要的一个方面，对个人、社会和经济发展都有深远的影响。它是一种工具，能让人们发挥自己的潜力，改善自己的生活质量，并为社会做出贡献。教育是每个人都应该拥有的基本权利，它对于创造一个公平和公正的社会是至关重要的。因此，我们投资于教育、确保每个人都能接受到优质的教育，这是至关重要的。", "总之，教育是人类生命中一个至关重要的方面，对个人、社会和经济发展都产生深远的影响。它是一个工具，使个人能够发挥潜力，改善生活质量并为社会做出贡献。教育是一项基本权利，应该让每个人都能够接受，这是创造一个公平和公正社会所必不可少的。因此，我们投资于教育并确保每个人都能够接受到优质的教育是至关重要的。"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}
